# Fix category modules and align income model with addIncome flow
#
# - Row 2 becomes a "bonus" income entry (was a leftover email-as-source
#   test row) with a refreshed amount + date.
# - A new "salary2" row is appended as row 3, mirroring the shape that the
#   addIncome flow writes (Source, Amount, Date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2: bonus / 2000 / 2026-02-17 ---
$ws.Range("A2").Value = "bonus"
$ws.Range("B2").Value = 2000
$ws.Range("C2").Value = 46070.250231481485

# --- Append new row 3: salary2 / 1000 / 2026-02-16 ---
$ws.Range("A3").Value = "salary2"
$ws.Range("B3").Value = 1000
$ws.Range("C3").Value = 46069.250231481485

# Give the new date cell (C3) the same date formatting/style as C2, the way
# addIncome would when it writes a new row under an existing formatted one.
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)
